$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 100.5
$ws.Range("I6").Value = 100.5
$ws.Range("K6").Value = 301.5
$ws.Range("M6").Value = -189.5
$ws.Range("H21").Value = 19642.846
$ws.Range("I21").Value = 24382.125
$ws.Range("J21").Value = 12060
$ws.Range("K21").Value = 24382.125
$ws.Range("L21").Value = 12060
$ws.Range("M21").Value = -23914.125
$ws.Range("N21").Value = -12996
$ws.Range("H23").Value = 19642.846
$ws.Range("I23").Value = 24382.125
$ws.Range("J23").Value = 12060
$ws.Range("K23").Value = 24382.125
$ws.Range("L23").Value = 12060
$ws.Range("M23").Value = -24148.125
$ws.Range("N23").Value = -12528
$ws.Range("H38").Value = 48.25
$ws.Range("I38").Value = 48.25
$ws.Range("K38").Value = 144.75
$ws.Range("M38").Value = 227.25
$ws.Range("H43").Value = 400
$ws.Range("I43").Value = 400
$ws.Range("J43").Value = 400
$ws.Range("K43").Value = 400
$ws.Range("L43").Value = 400
$ws.Range("M43").Value = -331
$ws.Range("N43").Value = -538
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H69").Value = 5780
$ws.Range("I69").Value = 4326.6665
$ws.Range("J69").Value = 6870
$ws.Range("K69").Value = 12979.9995
$ws.Range("L69").Value = 20610
$ws.Range("M69").Value = -12105.9995
$ws.Range("N69").Value = -22358
$ws.Range("H72").Value = 5780
$ws.Range("I72").Value = 4326.6665
$ws.Range("J72").Value = 6870
$ws.Range("K72").Value = 38939.9985
$ws.Range("L72").Value = 61830
$ws.Range("M72").Value = -34571.9985
$ws.Range("N72").Value = -70566
$ws.Range("H112").Value = 41668030
$ws.Range("I112").Value = 699.8
$ws.Range("J112").Value = 52633120
$ws.Range("K112").Value = 2099.4
$ws.Range("L112").Value = 157899360
$ws.Range("M112").Value = -991.3999999999996
$ws.Range("N112").Value = -157901576
$ws.Range("H126").Value = 60780
$ws.Range("J126").Value = 60780
$ws.Range("L126").Value = 60780
$ws.Range("N126").Value = -70660
$ws.Range("H129").Value = 1128.75
$ws.Range("I129").Value = 302.54544
$ws.Range("J129").Value = 10217
$ws.Range("K129").Value = 907.63632
$ws.Range("L129").Value = 30651
$ws.Range("M129").Value = 4092.36368
$ws.Range("N129").Value = -40651
$ws.Range("H137").Value = 30306698
$ws.Range("I137").Value = 1967.0667
$ws.Range("J137").Value = 55560640
$ws.Range("K137").Value = 5901.2001
$ws.Range("L137").Value = 166681920
$ws.Range("M137").Value = -3351.2001
$ws.Range("N137").Value = -166687020
$ws.Range("H138").Value = 3324.7446
$ws.Range("I138").Value = 3419.8928
$ws.Range("J138").Value = 3284.379
$ws.Range("K138").Value = 10259.6784
$ws.Range("L138").Value = 9853.136999999999
$ws.Range("M138").Value = -5119.678400000001
$ws.Range("N138").Value = -20133.137
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 12027.625
$ws.Range("I28").Value = 6602.857
$ws.Range("K28").Value = 6602.857
$ws.Range("M28").Value = -6410.857
$ws.Range("H32").Value = 19833.658
$ws.Range("I32").Value = 14346.409
$ws.Range("J32").Value = 71570.57000000001
$ws.Range("K32").Value = 14346.409
$ws.Range("L32").Value = 71570.57000000001
$ws.Range("M32").Value = -14059.409
$ws.Range("N32").Value = -72144.57000000001
$ws.Range("H99").Value = 12027.625
$ws.Range("I99").Value = 6602.857
$ws.Range("K99").Value = 6602.857
$ws.Range("M99").Value = -3607.857
$ws.Range("H132").Value = 842677.4399999999
$ws.Range("I132").Value = 986787.7
$ws.Range("J132").Value = 2034.3334
$ws.Range("K132").Value = 2960363.1
$ws.Range("L132").Value = 6103.0002
$ws.Range("M132").Value = -2957833.1
$ws.Range("N132").Value = -11163.0002
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 32090.5
$ws.Range("J62").Value = 32090.5
$ws.Range("L62").Value = 32090.5
$ws.Range("N62").Value = -33462.5
$ws.Range("H65").Value = 32090.5
$ws.Range("J65").Value = 32090.5
$ws.Range("L65").Value = 96271.5
$ws.Range("N65").Value = -103135.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 2000
$ws.Range("J17").Value = 2000
$ws.Range("L17").Value = 2000
$ws.Range("N17").Value = -2348
$ws.Range("H22").Value = 205.55556
$ws.Range("I22").Value = 198.69565
$ws.Range("J22").Value = 245
$ws.Range("K22").Value = 198.69565
$ws.Range("L22").Value = 245
$ws.Range("M22").Value = 151.30435
$ws.Range("N22").Value = -945
$ws.Range("H25").Value = 2337.6667
$ws.Range("J25").Value = 6013
$ws.Range("L25").Value = 6013
$ws.Range("N25").Value = -6361
$ws.Range("H50").Value = 11219
$ws.Range("J50").Value = 11219
$ws.Range("L50").Value = 11219
$ws.Range("N50").Value = -12469
$ws.Range("H74").Value = 14421.2
$ws.Range("J74").Value = 16730.25
$ws.Range("L74").Value = 16730.25
$ws.Range("N74").Value = -18478.25
$ws.Range("H77").Value = 14421.2
$ws.Range("J77").Value = 16730.25
$ws.Range("L77").Value = 50190.75
$ws.Range("N77").Value = -58926.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 2000
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 2000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 6000
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -6460
$ws.Range("H68").Value = 768.5
$ws.Range("I68").Value = 609.4286
$ws.Range("J68").Value = 927.5714
$ws.Range("K68").Value = 1828.2858
$ws.Range("L68").Value = 2782.7142
$ws.Range("M68").Value = -1017.2858
$ws.Range("N68").Value = -4404.7142
$ws.Range("H71").Value = 768.5
$ws.Range("I71").Value = 609.4286
$ws.Range("J71").Value = 927.5714
$ws.Range("K71").Value = 5484.8574
$ws.Range("L71").Value = 8348.142600000001
$ws.Range("M71").Value = -1428.8574
$ws.Range("N71").Value = -16460.1426
$ws.Range("H102").Value = 750
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H121").Value = 27778634
$ws.Range("I121").Value = 616.25
$ws.Range("J121").Value = 50001050
$ws.Range("K121").Value = 1848.75
$ws.Range("L121").Value = 150003150
$ws.Range("M121").Value = -538.75
$ws.Range("N121").Value = -150005770
$ws.Range("H132").Value = 55556788
$ws.Range("I132").Value = 76924136
$ws.Range("J132").Value = 1682
$ws.Range("K132").Value = 692317224
$ws.Range("L132").Value = 15138
$ws.Range("M132").Value = -692314694
$ws.Range("N132").Value = -20198
$ws.Range("H137").Value = 25384.756
$ws.Range("I137").Value = 4048.3333
$ws.Range("J137").Value = 28361.93
$ws.Range("K137").Value = 12144.9999
$ws.Range("L137").Value = 85085.79000000001
$ws.Range("M137").Value = -7044.999899999999
$ws.Range("N137").Value = -95285.79000000001
$ws.Range("H138").Value = 3702.2068
$ws.Range("I138").Value = 2290.3572
$ws.Range("J138").Value = 5019.933
$ws.Range("K138").Value = 6871.071599999999
$ws.Range("L138").Value = 15059.799
$ws.Range("M138").Value = -1731.071599999999
$ws.Range("N138").Value = -25339.799
$ws.Range("H139").Value = 5266.6665
$ws.Range("I139").Value = 5266.6665
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 15799.9995
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -10659.9995
$ws.Range("N139").ClearContents()
